# Applies the "200524 23:13" commit: adds a new "loginTestInValid" row to the
# Runner sheet, adds a "Browser" column to the CredentialData sheet, and
# updates a handful of existing values/selections/views to match.

$wb = $excel.ActiveWorkbook
$wsRunner = $wb.Worksheets.Item("Runner")
$wsCred   = $wb.Worksheets.Item("CredentialData")

# --- Runner sheet: new row 4 ---------------------------------------------
$wsRunner.Range("A4").Value = "loginTestInValid"
$wsRunner.Range("B4").Value = "To verify login functionality with invalid credentials"
$wsRunner.Range("C4").Value = "Yes"
$wsRunner.Range("D4").Value = "'1"
$wsRunner.Range("E4").Value = "'1"

# --- CredentialData sheet: fix Execute column + new Browser column -------
$wsCred.Range("D2").Value = "Yes"
$wsCred.Range("D3").Value = "Yes"

$wsCred.Range("E1").Value = "Browser"
$wsCred.Range("E2").Value = "chrome"
$wsCred.Range("E3").Value = "firefox"
$wsCred.Range("E4").Value = "edge"

# --- Selections / views ---------------------------------------------------
$wsRunner.Range("E10").Select()
$wsCred.Range("H7").Select()
